$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 13 (shifts rows 13-23 down to 14-24), which is where
# the "Docentes responsáveis:" value (the professor's name, previously
# misplaced in row 10 alongside "Objetivos:") now lives.
$ws.Rows(13).Insert()

# Row 10 ("Objetivos:") gets its own real objective text instead of the
# professor name that had been placed there by mistake.
$ws.Range("B10").Value = 'Desenvolver o aprendizado teórico e prático da Bioquímica através da execução de práticas de laboratório baseadas na evolução do conteúdo teórico ministrado na disciplina Bioquímica II.'
$ws.Range("C10").Value = 'Desenvolver o aprendizado teórico e prático da Bioquímica através da execução de práticas de laboratório baseadas na evolução do conteúdo teórico ministrado na disciplina Bioquímica II.'

# Newly-inserted row 13, under "Docentes responsáveis:" (row 12), now holds
# the professor's name/ID.
$ws.Range("B13").Value = '6007846 - Júlio César dos Santos'
$ws.Range("C13").Value = '6007846 - Júlio César dos Santos'

# Row 14 ("Programa resumido:") gets the real short-syllabus text instead of
# the placeholder "Semestral".
$ws.Range("B14").Value = 'Propriedades gerais de glicídios; Fermentação anaeróbia; Extração deClorofila e Reação de Hill; Transporte de glicídios e indução de enzimas.'
$ws.Range("C14").Value = 'Propriedades gerais de glicídios; Fermentação anaeróbia; Extração deClorofila e Reação de Hill; Transporte de glicídios e indução de enzimas.'

# Row 16 ("Programa:") gets the real full syllabus text instead of the
# placeholder date "01/01/2018".
$ws.Range("B16").Value = 'Propriedades gerais de glicídios:principais testes qualitativos para identificação e diferenciação de glicídios; aplicação de certas reações coloridas e dosagem espectrofotométrica de monossacarídeos redutores.Fermentação anaeróbia: Conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono;cálculo da eficiência do processo; ação de um inibidor da glicólise. Extração de clorofila e reação deHill: estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura e fase luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. Transporte de glicídios e indução de enzimas: conceitos gerais; enzimas do catabolismo da galactose; repressão, inativação emodificação catabólicas; sistemas enzimáticos constitutivos e induzidos em células de levedura'
$ws.Range("C16").Value = 'Propriedades gerais de glicídios:principais testes qualitativos para identificação e diferenciação de glicídios; aplicação de certas reações coloridas e dosagem espectrofotométrica de monossacarídeos redutores.Fermentação anaeróbia: Conceitos gerais e fermentação de glicose; produção de etanol e CO2; consumo da fonte de carbono;cálculo da eficiência do processo; ação de um inibidor da glicólise. Extração de clorofila e reação deHill: estrutura de cloroplastos; papel da clorofila nos sistemas fotossintéticos I e II; fase escura e fase luminosa; produção de NADP; produção de ATP; papel do corante como aceptor de prótons e elétrons. Transporte de glicídios e indução de enzimas: conceitos gerais; enzimas do catabolismo da galactose; repressão, inativação emodificação catabólicas; sistemas enzimáticos constitutivos e induzidos em células de levedura'

# Row 19 ("Método:") gets the evaluation-method text that used to sit one
# row down, next to "Critério:".
$ws.Range("B19").Value = 'A avaliação será feita por meio de uma prova escrita e notas de relatórios (R).'
$ws.Range("C19").Value = 'A avaliação será feita por meio de uma prova escrita e notas de relatórios (R).'

# Row 20 ("Critério:") gets the final-grade formula text that used to sit
# one row down, next to "Norma de recuperação:".
$ws.Range("B20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1*2 + R)/3.'
$ws.Range("C20").Value = 'A Nota final (NF) será calculada da seguinte maneira: NF = (P1*2 + R)/3.'

# Row 21 ("Norma de recuperação:") gets the makeup-exam text that used to
# sit one row down, next to "Bibliografia:".
$ws.Range("B21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada pela fórmula: MR = (NF + PR)/2'
$ws.Range("C21").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) será calculada pela fórmula: MR = (NF + PR)/2'

# Row 22 ("Bibliografia:") gets the real bibliography text.
$ws.Range("B22").Value = 'CISTERNAS, J. R. Fundamentos de bioquímica experimental. São Paulo : Atheneu, 2005. ISBN: 9788573791075.NELSON, D. L., COX. M. M. Princípios de bioquímica de Lehninger. Porto Alegre : Artmed, 2011. ISBN: 9788536324180.VOET, D., VOET, J. G. Bioquímica. Porto Alegre : Artmed, 2013. ISBN: 9788582710043.'
$ws.Range("C22").Value = 'CISTERNAS, J. R. Fundamentos de bioquímica experimental. São Paulo : Atheneu, 2005. ISBN: 9788573791075.NELSON, D. L., COX. M. M. Princípios de bioquímica de Lehninger. Porto Alegre : Artmed, 2011. ISBN: 9788536324180.VOET, D., VOET, J. G. Bioquímica. Porto Alegre : Artmed, 2013. ISBN: 9788582710043.'
